# Automatische test-sync: 2025-06-24 22:05:50
# Append a new log row (row 41) to the "Logs" sheet and update the
# "Dashboard" summary count for the "Factuur / Administratie" category.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Add the new row of data to the Logs sheet ---
$logs.Range("A41").Value = "Herinnering betaling"
$logs.Range("B41").Value = "mailmind.test@zohomail.eu"
$logs.Range("C41").Value = "Ik zie dat ik nog een openstaande betaling heb. Kunt u dit bevestigen?"
$logs.Range("D41").Value = "Factuur / Administratie"
$logs.Range("F41").Value = "2025-06-24 22:05:29"
$logs.Range("G41").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$catFormats = $logs.Range("D2:D40").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D41"))
}

$answeredFormats = $logs.Range("G2:G40").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G41"))
}

# --- Update the Dashboard count for "Factuur / Administratie" ---
$dashboard.Range("B3").Value = 7
